$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7554823756217957
$ws.Range("B1").Value = 3.281688928604126
$ws.Range("C1").Value = 2.95831298828125
$ws.Range("D1").Value = 2.357467412948608
$ws.Range("E1").Value = 1.477815866470337
